# Simulated Wild Card round and logged it
# Update row 2 (Houston "H") stats on both OFF and DEF sheets
# with the results of the simulated Wild Card round game.

$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 234
$wsOff.Range("C2").Value = 166
$wsOff.Range("D2").Value = 60
$wsOff.Range("E2").Value = 33

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 196
$wsDef.Range("C2").Value = 138
$wsDef.Range("D2").Value = 47
$wsDef.Range("E2").Value = 30

$wb.Save()
